# Katalog.xlsx update:
#  - Build out the "GeländerEdelstahl_Daten" pricing-catalog sheet (Edelstahl
#    Geländer / stainless-steel railing configurator rows), mirroring the
#    layout already used on "Carport_Liste".
#  - Make "GeländerEdelstahl_Daten" the active/selected sheet & tab.
#  - Leave the previous active sheet ("Carport_Liste") selection pointed at
#    its full used range.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Carport_Liste")
$ws4 = $wb.Worksheets.Item("GeländerEdelstahl_Daten")

# ---------------------------------------------------------------------
# Style "donor" cells already present on Carport_Liste - one per cellXfs
# index we need (1, 2, 4, 5). Copying formats from these (PasteSpecial
# xlPasteFormats = -4122) reproduces the exact same style index on the
# destination cell without touching its value.
# ---------------------------------------------------------------------
$styleDonor = @{
    1 = $ws2.Range("B2")
    2 = $ws2.Range("A1")
    4 = $ws2.Range("C2")
    5 = $ws2.Range("E2")
}

function Set-CellStyle($range, $styleId) {
    $styleDonor[$styleId].Copy()
    $range.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Column widths (B/D/E get wider to fit the new option/formula text).
# ---------------------------------------------------------------------
$ws4.Columns(2).ColumnWidth = 30.5
$ws4.Columns(4).ColumnWidth = 32
$ws4.Columns(5).ColumnWidth = 42.166666666666664

# ---------------------------------------------------------------------
# Row data: column -> (value, styleId). $null value = leave cell blank
# (style only).
# ---------------------------------------------------------------------
$rows = @{
    1  = @{ A=@("Typ",2);       B=@("Bezeichnung",2); C=@("Variable",2); D=@("Optionen",2); E=@("O_Faktor",2) }
    2  = @{ A=@("Zahl",2);      B=@("Länge des Geländers (lfm)",1); C=@("L",4);      D=@($null,1); E=@($null,1) }
    3  = @{ A=@("Zahl",2);      B=@("Anzahl Ecken (Stk)",1);        C=@("Ecken",4);  D=@($null,1); E=@($null,1) }
    4  = @{ A=@("Auswahl",2);   B=@("Profil-Type (Steher/Handlauf)",1); C=@("P_Basis",4); D=@("Rund 42.4mm:100, Rund 48.3mm:120, Quadrat 40x40:110, Quadrat 50x50:140",4); E=@("Basispreis für Steher/Handlauf pro Meter",5) }
    5  = @{ A=@("Auswahl",2);   B=@("Füllung",1);                   C=@("P_Full",4); D=@("Stäbe Horizontal (OIB prüfen!):40, Stäbe Vertikal:120, Glas VSG Klar:180, Glas VSG Matt:200, Seile:30, Lochblech:140",4); E=@("Preis der Füllung pro Meter",5) }
    6  = @{ A=@("Auswahl",2);   B=@("Montageart",1);                C=@("P_Mount",4); D=@("Von Oben (Boden):0, Von Seite (Stirn):25",4); E=@("Aufpreis pro Meter für Seitenmontage",5) }
    7  = @{ A=@("Auswahl",2);   B=@("Deckenblende / Verblendung",1); C=@("P_Blende",4); D=@("Keine:0, Alu-Blende:50, Edelstahl-Blende:90",4); E=@("Preis pro Meter für Rand-Abdeckung",5) }
    8  = @{ A=@("Auswahl",2);   B=@("Montage gewünscht?",1);        C=@("M_Faktor",4); D=@("Ja:1, Nein:0",4); E=@($null,1) }
    9  = @{ A=@("Zahl",2);      B=@("Montagepreis (€/lfm)",1);      C=@("M_Preis",4); D=@($null,1); E=@("Dein Satz für Montage pro Meter",5) }
    10 = @{ A=@("Preis",2);     B=@("Gesamtpreis",2);               C=@("Formel / Info",4); D=@($null,1); E=@("(L * (P_Basis + P_Full + P_Mount + P_Blende)) + (Ecken * 100) + (L * M_Preis * M_Faktor)",2) }
    11 = @{ A=@($null,2);       B=@($null,1);                       C=@($null,4); D=@($null,4); E=@($null,5) }
    12 = @{ A=@($null,2);       B=@($null,1);                       C=@($null,4); D=@($null,1); E=@($null,5) }
    13 = @{ A=@($null,2);       B=@($null,2);                       C=@($null,4); D=@($null,1); E=@($null,2) }
}

$rowHeights = @{
    1=30.75; 2=43.5; 3=29.25; 4=115.5; 5=179.25; 6=51.75; 7=64.5; 8=29.25; 9=29.25; 10=15.75; 11=15.75; 12=15.75; 13=15.75
}

$cols = @("A","B","C","D","E")

for ($r = 1; $r -le 13; $r++) {
    $rowData = $rows[$r]
    foreach ($col in $cols) {
        $pair = $rowData[$col]
        $value = $pair[0]
        $styleId = $pair[1]
        $cellRange = $ws4.Range($col + $r)
        Set-CellStyle $cellRange $styleId
        if ($value -ne $null) {
            $cellRange.Value = $value
        }
    }
    $ws4.Rows($r).RowHeight = $rowHeights[$r]
}

# ---------------------------------------------------------------------
# Sheet/window state: Carport_Liste's old selection clears, the whole
# used range A1:E13 becomes selected instead; GeländerEdelstahl_Daten
# becomes the active tab with E11 selected.
# ---------------------------------------------------------------------
$ws2.Range("A1:E13").Select()
$ws4.Activate()
$ws4.Range("E11").Select()

Write-Host "GeländerEdelstahl_Daten populated."
